# Feria Lagunitas de Puerto Montt - Pepino ensalada
# Insert a new weekly data row right before the current row 130, shifting
# the existing rows 130-142 down to 131-143 (this mirrors the diff, where
# every later row's content becomes identical to what used to be in the
# row right above it, and a brand-new row is introduced at the old row 130
# position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 130 and below down by one row.
$ws.Rows("130:130").Insert()

# Populate the newly inserted row 130 with the new weekly record.
$ws.Cells.Item(130, 1).Value = 4
$ws.Cells.Item(130, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(130, 3).Value = "Los Lagos"
$ws.Cells.Item(130, 4).Value = 44449
$ws.Cells.Item(130, 5).Value = 10
$ws.Cells.Item(130, 6).Value = 100112043
$ws.Cells.Item(130, 7).Value = "Pepino ensalada"
$ws.Cells.Item(130, 8).Value = "Sin especificar"
$ws.Cells.Item(130, 9).Value = "Primera"
$ws.Cells.Item(130, 10).Value = 300
$ws.Cells.Item(130, 11).Value = 21000
$ws.Cells.Item(130, 12).Value = 21000
$ws.Cells.Item(130, 13).Value = 21000
$ws.Cells.Item(130, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(130, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(130, 16).Value = 350
$ws.Cells.Item(130, 17).Value = 60
$ws.Cells.Item(130, 18).Value = "Hortaliza"

# Match the date cell style used by the rest of column D (numFmtId 165).
$ws.Cells.Item(130, 4).NumberFormat = $ws.Cells.Item(131, 4).NumberFormat
